$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 1172626.4
$ws.Range("I64").Value = 2108143
$ws.Range("J64").Value = 3230.5
$ws.Range("K64").Value = 2108143
$ws.Range("L64").Value = 3230.5
$ws.Range("M64").Value = -2107895
$ws.Range("N64").Value = -3726.5

$ws.Range("H67").Value = 1172626.4
$ws.Range("I67").Value = 2108143
$ws.Range("J67").Value = 3230.5
$ws.Range("K67").Value = 2108143
$ws.Range("L67").Value = 3230.5
$ws.Range("M67").Value = -2107285
$ws.Range("N67").Value = -4946.5

$ws.Range("H76").Value = 3370668
$ws.Range("I76").Value = 3707185
$ws.Range("J76").Value = 5500
$ws.Range("K76").Value = 3707185
$ws.Range("L76").Value = 5500
$ws.Range("M76").Value = -3706870
$ws.Range("N76").Value = -6130

$ws.Range("H79").Value = 3370668
$ws.Range("I79").Value = 3707185
$ws.Range("J79").Value = 5500
$ws.Range("K79").Value = 3707185
$ws.Range("L79").Value = 5500
$ws.Range("M79").Value = -3706093
$ws.Range("N79").Value = -7684

$ws.Range("H128").Value = 32840
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 32840
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 32840
$ws.Range("N128").Value = -42800

$ws.Range("H138").Value = 2268.027
$ws.Range("I138").Value = 1714.5667
$ws.Range("J138").Value = 4640
$ws.Range("K138").Value = 5143.7001
$ws.Range("L138").Value = 13920
$ws.Range("M138").Value = -3.70010000000002
$ws.Range("N138").Value = -24200

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 337793.47
$ws.Range("I32").Value = 4391.8237
$ws.Range("J32").Value = 1597310.8
$ws.Range("K32").Value = 4391.8237
$ws.Range("L32").Value = 1597310.8
$ws.Range("M32").Value = -4104.8237
$ws.Range("N32").Value = -1597884.8

$ws.Range("H61").Value = 2539.225
$ws.Range("I61").Value = 2690.3794
$ws.Range("J61").Value = 2140.7273
$ws.Range("K61").Value = 2690.3794
$ws.Range("L61").Value = 2140.7273
$ws.Range("M61").Value = -2478.3794
$ws.Range("N61").Value = -2564.7273

$ws.Range("H74").Value = 1336.6666
$ws.Range("I74").Value = 938.5217
$ws.Range("J74").Value = 2644.8572
$ws.Range("K74").Value = 938.5217
$ws.Range("L74").Value = 2644.8572
$ws.Range("M74").Value = -64.52170000000001
$ws.Range("N74").Value = -4392.8572

$ws.Range("H77").Value = 1336.6666
$ws.Range("I77").Value = 938.5217
$ws.Range("J77").Value = 2644.8572
$ws.Range("K77").Value = 4692.6085
$ws.Range("L77").Value = 13224.286
$ws.Range("M77").Value = -324.6085000000003
$ws.Range("N77").Value = -21960.286

$ws.Range("H88").Value = 9245
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 9245
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 9245
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -10057

$ws.Range("H91").Value = 9245
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 9245
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 9245
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -12053

$ws.Range("H109").Value = 44994.184
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 44994.184
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 44994.184
$ws.Range("N109").Value = -47768.184

$ws.Range("H132").Value = 1942.2778
$ws.Range("I132").Value = 1526.7693
$ws.Range("J132").Value = 3022.6
$ws.Range("K132").Value = 4580.3079
$ws.Range("L132").Value = 9067.799999999999
$ws.Range("M132").Value = -2050.3079
$ws.Range("N132").Value = -14127.8

$ws.Range("H136").Value = 2539.225
$ws.Range("I136").Value = 2690.3794
$ws.Range("J136").Value = 2140.7273
$ws.Range("K136").Value = 8071.138199999999
$ws.Range("L136").Value = 6422.1819
$ws.Range("M136").Value = -5521.138199999999
$ws.Range("N136").Value = -11522.1819

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3271.24
$ws.Range("I62").Value = 3369.5264
$ws.Range("J62").Value = 2960
$ws.Range("K62").Value = 3369.5264
$ws.Range("L62").Value = 2960
$ws.Range("M62").Value = -2745.5264
$ws.Range("N62").Value = -4208

$ws.Range("H65").Value = 3271.24
$ws.Range("I65").Value = 3369.5264
$ws.Range("J65").Value = 2960
$ws.Range("K65").Value = 16847.632
$ws.Range("L65").Value = 14800
$ws.Range("M65").Value = -13727.632
$ws.Range("N65").Value = -21040

$ws.Range("H99").Value = 1816.0834
$ws.Range("I99").Value = 1750.849
$ws.Range("J99").Value = 2310
$ws.Range("K99").Value = 1750.849
$ws.Range("L99").Value = 2310
$ws.Range("M99").Value = -252.8489999999999
$ws.Range("N99").Value = -5306

$ws.Range("H126").Value = 1816.0834
$ws.Range("I126").Value = 1750.849
$ws.Range("J126").Value = 2310
$ws.Range("K126").Value = 5252.547
$ws.Range("L126").Value = 6930
$ws.Range("M126").Value = -2782.547
$ws.Range("N126").Value = -11870

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 15209540
$ws.Range("I70").Value = 22505852
$ws.Range("J70").Value = 8891.583000000001
$ws.Range("K70").Value = 22505852
$ws.Range("L70").Value = 8891.583000000001
$ws.Range("M70").Value = -22505582
$ws.Range("N70").Value = -9431.583000000001

$ws.Range("H73").Value = 15209540
$ws.Range("I73").Value = 22505852
$ws.Range("J73").Value = 8891.583000000001
$ws.Range("K73").Value = 22505852
$ws.Range("L73").Value = 8891.583000000001
$ws.Range("M73").Value = -22504916
$ws.Range("N73").Value = -10763.583

$ws.Range("H80").Value = 3500
$ws.Range("I80").Value = 2000
$ws.Range("J80").Value = 4000
$ws.Range("K80").Value = 2000
$ws.Range("L80").Value = 4000
$ws.Range("M80").Value = -1002
$ws.Range("N80").Value = -5996

$ws.Range("H83").Value = 3500
$ws.Range("I83").Value = 2000
$ws.Range("J83").Value = 4000
$ws.Range("K83").Value = 10000
$ws.Range("L83").Value = 20000
$ws.Range("M83").Value = -5008
$ws.Range("N83").Value = -29984

$ws.Range("H132").Value = 2089.862
$ws.Range("I132").Value = 1399.7222
$ws.Range("J132").Value = 3219.182
$ws.Range("K132").Value = 4199.1666
$ws.Range("L132").Value = 9657.545999999998
$ws.Range("M132").Value = -1669.1666
$ws.Range("N132").Value = -14717.546

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1755.1143
$ws.Range("I40").Value = 1691.15
$ws.Range("J40").Value = 1840.4
$ws.Range("K40").Value = 1691.15
$ws.Range("L40").Value = 1840.4
$ws.Range("M40").Value = -1555.15
$ws.Range("N40").Value = -2112.4

$ws.Range("H122").Value = 2264.3809
$ws.Range("I122").Value = 1986.1765
$ws.Range("J122").Value = 3446.75
$ws.Range("K122").Value = 5958.529500000001
$ws.Range("L122").Value = 10340.25
$ws.Range("M122").Value = -3508.529500000001
$ws.Range("N122").Value = -15240.25

$ws.Range("H136").Value = 5264.3335
$ws.Range("I136").Value = 3347.0908
$ws.Range("J136").Value = 8277.143
$ws.Range("K136").Value = 10041.2724
$ws.Range("L136").Value = 24831.429
$ws.Range("M136").Value = -7491.2724
$ws.Range("N136").Value = -29931.429

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 13437.429
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 13437.429
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 13437.429
$ws.Range("M45").ClearContents()
$ws.Range("N45").Value = -14419.429

$ws.Range("H80").Value = 30000
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 30000
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 30000
$ws.Range("N80").Value = -31996

$ws.Range("H82").Value = 10000
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 10000
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 10000
$ws.Range("N82").Value = -10766

$ws.Range("H83").Value = 30000
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 30000
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 90000
$ws.Range("N83").Value = -99984

$ws.Range("H85").Value = 10000
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 10000
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 10000
$ws.Range("N85").Value = -12652

$ws.Range("H107").Value = 664.4
$ws.Range("I107").Value = 748.6
$ws.Range("J107").Value = 496
$ws.Range("K107").Value = 2245.8
$ws.Range("L107").Value = 1488
$ws.Range("M107").Value = -325.8000000000002
$ws.Range("N107").Value = -5328

$ws.Range("H122").Value = 1466.6923
$ws.Range("I122").Value = 1568.7
$ws.Range("J122").Value = 1126.6666
$ws.Range("K122").Value = 4706.1
$ws.Range("L122").Value = 3379.9998
$ws.Range("M122").Value = -2256.1
$ws.Range("N122").Value = -8279.9998

$ws.Range("H126").Value = 1366.2572
$ws.Range("I126").Value = 1185.4333
$ws.Range("J126").Value = 2451.2
$ws.Range("K126").Value = 3556.2999
$ws.Range("L126").Value = 7353.599999999999
$ws.Range("M126").Value = -1086.2999
$ws.Range("N126").Value = -12293.6

$ws.Range("H132").Value = 50005180
$ws.Range("I132").Value = 68183384
$ws.Range("J132").Value = 15126
$ws.Range("K132").Value = 204550152
$ws.Range("L132").Value = 45378
$ws.Range("M132").Value = -204547622
$ws.Range("N132").Value = -50438

$ws.Range("H136").Value = 741.84906
$ws.Range("I136").Value = 524.0244
$ws.Range("J136").Value = 1486.0834
$ws.Range("K136").Value = 1572.0732
$ws.Range("L136").Value = 4458.2502
$ws.Range("M136").Value = 977.9268
$ws.Range("N136").Value = -9558.2502
